$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 203 is a brand-new row at the end of the table; clone row 202 (values + formatting)
# into it first so styles (e.g. the date format on column D) carry over correctly.
$ws.Range("A202:R202").Copy($ws.Range("A203:R203"))

# Columns D and I:P hold the per-record data (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad, Origen, Precio $/Kg). Each record from
# row 103 downward is now one row lower, and row 103 itself gets a brand-new Fecha value,
# so rewrite that block for rows 103-203 with the final values.
$rows = @(
    @(44586, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44252, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44356, "Primera", 10, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44356, "Primera", 20, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833),
    @(44487, "Primera", 105, 5000, 6000, 5619, "`$/docena de atados", "Provincia de Cautín", 1873),
    @(44487, "Primera", 75, 2000, 2000, 2000, "`$/docena de atados", "Región Metropolitana", 667),
    @(44410, "Primera", 30, 8000, 8000, 8000, "`$/docena de atados", "Provincia de Cautín", 2667),
    @(44327, "Primera", 20, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44455, "Primera", 20, 7000, 8000, 7500, "`$/docena de atados", "Provincia de Cautín", 2500),
    @(44582, "Primera", 50, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44159, "Primera", 40, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44235, "Primera", 110, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44461, "Primera", 40, 3000, 6000, 3750, "`$/docena de atados", "Provincia de Cautín", 1250),
    @(44466, "Primera", 40, 7000, 7000, 7000, "`$/docena de atados", "Provincia de Cautín", 2333),
    @(44466, "Primera", 50, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000),
    @(44462, "Primera", 40, 6000, 7000, 6500, "`$/docena de atados", "Provincia de Cautín", 2167),
    @(44580, "Primera", 40, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44364, "Primera", 65, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833),
    @(44463, "Primera", 30, 6000, 7000, 6667, "`$/docena de atados", "Provincia de Cautín", 2222),
    @(44309, "Primera", 55, 7000, 7000, 7000, "`$/docena de atados", "Provincia de Cautín", 2333),
    @(44484, "Primera", 35, 5000, 6000, 5571, "`$/docena de atados", "Provincia de Cautín", 1857),
    @(44313, "Primera", 30, 5000, 6000, 5333, "`$/docena de atados", "Provincia de Cautín", 1778),
    @(44516, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44452, "Primera", 30, 8000, 8000, 8000, "`$/docena de atados", "Provincia de Cautín", 2667),
    @(44452, "Segunda", 20, 5000, 5000, 5000, "`$/docena de atados", "Región Metropolitana", 1667),
    @(44257, "Primera", 30, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44322, "Primera", 65, 5000, 6000, 5462, "`$/docena de atados", "Provincia de Cautín", 1821),
    @(44246, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44218, "Primera", 50, 4000, 5000, 4400, "`$/docena de atados", "Provincia de Cautín", 1467),
    @(44400, "Primera", 10, 10000, 10000, 10000, "`$/docena de atados", "Provincia de Cautín", 3333),
    @(44426, "Primera", 10, 8000, 8000, 8000, "`$/docena de atados", "Provincia de Cautín", 2667),
    @(44431, "Primera", 55, 8000, 8000, 8000, "`$/docena de atados", "Provincia de Cautín", 2667),
    @(44211, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44217, "Primera", 130, 4000, 5000, 4385, "`$/docena de atados", "Provincia de Cautín", 1462),
    @(44512, "Primera", 20, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44511, "Primera", 50, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44193, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44166, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44200, "Primera", 20, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Cautín", 1333),
    @(44175, "Primera", 40, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Cautín", 1333),
    @(44160, "Primera", 20, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44203, "Primera", 80, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Cautín", 1333),
    @(44519, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44161, "Primera", 30, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44386, "Primera", 10, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44253, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44581, "Primera", 150, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44504, "Primera", 125, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44567, "Primera", 40, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44280, "Primera", 95, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Cautín", 1333),
    @(44442, "Primera", 15, 10000, 10000, 10000, "`$/docena de atados", "Provincia de Cautín", 3333),
    @(44344, "Primera", 30, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44536, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44187, "Primera", 30, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Cautín", 1333),
    @(44277, "Primera", 65, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Cautín", 1333),
    @(44391, "Primera", 55, 7000, 7000, 7000, "`$/docena de atados", "Provincia de Cautín", 2333),
    @(44454, "Primera", 20, 8000, 8000, 8000, "`$/docena de atados", "Provincia de Cautín", 2667),
    @(44202, "Primera", 20, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Cautín", 1333),
    @(44371, "Primera", 40, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44249, "Primera", 110, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44259, "Primera", 80, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44328, "Primera", 45, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44396, "Primera", 30, 7000, 7000, 7000, "`$/docena de atados", "Provincia de Cautín", 2333),
    @(44494, "Primera", 20, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44494, "Primera", 30, 2000, 2000, 2000, "`$/docena de atados", "Región Metropolitana", 667),
    @(44526, "Primera", 20, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44264, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44533, "Primera", 65, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44354, "Primera", 30, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44221, "Primera", 85, 5000, 6000, 5588, "`$/docena de atados", "Provincia de Cautín", 1863),
    @(44523, "Primera", 30, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44399, "Primera", 30, 10000, 10000, 10000, "`$/docena de atados", "Provincia de Cautín", 3333),
    @(44382, "Primera", 30, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44441, "Primera", 20, 10000, 10000, 10000, "`$/docena de atados", "Provincia de Cautín", 3333),
    @(44167, "Primera", 55, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44372, "Primera", 20, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44335, "Primera", 35, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44316, "Primera", 20, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44475, "Primera", 30, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833),
    @(44315, "Primera", 40, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44186, "Primera", 40, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Cautín", 1333),
    @(44176, "Primera", 70, 4000, 5000, 4429, "`$/docena de atados", "Provincia de Cautín", 1476),
    @(44278, "Primera", 35, 400, 400, 400, "`$/docena de atados", "Provincia de Cautín", 133),
    @(44438, "Primera", 20, 10000, 10000, 10000, "`$/docena de atados", "Provincia de Cautín", 3333),
    @(44312, "Primera", 20, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44260, "Primera", 30, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44449, "Primera", 95, 8000, 8000, 8000, "`$/docena de atados", "Provincia de Cautín", 2667),
    @(44449, "Primera", 85, 4000, 4000, 4000, "`$/docena de atados", "Región Metropolitana", 1333),
    @(44474, "Primera", 30, 5000, 5000, 5000, "`$/docena de atados", "Región Metropolitana", 1667),
    @(44585, "Primera", 35, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44385, "Primera", 30, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44162, "Primera", 40, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44529, "Primera", 110, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44323, "Primera", 50, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44189, "Primera", 40, 4000, 5000, 4500, "`$/docena de atados", "Provincia de Cautín", 1500),
    @(44428, "Primera", 10, 8000, 8000, 8000, "`$/docena de atados", "Provincia de Cautín", 2667),
    @(44321, "Primera", 45, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000),
    @(44302, "Primera", 20, 7000, 7000, 7000, "`$/docena de atados", "Provincia de Cautín", 2333),
    @(44209, "Primera", 50, 8000, 8000, 8000, "`$/docena de atados", "Provincia de Cautín", 2667),
    @(44274, "Primera", 20, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667),
    @(44392, "Primera", 65, 3000, 3000, 3000, "`$/docena de atados", "Provincia de Cautín", 1000)
)

$startRow = 103
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 4).Value = $row[0]   # D Fecha
    $ws.Cells.Item($r, 9).Value = $row[1]   # I Calidad
    $ws.Cells.Item($r, 10).Value = $row[2]  # J Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]  # K Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]  # L Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]  # M Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $row[6]  # N Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $row[7]  # O Origen
    $ws.Cells.Item($r, 16).Value = $row[8]  # P Precio $/Kg
}

